$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Column F header text changes from "variance" to "confidence interval" wording.
$ws.Range("F1").Value = "Доверительный интервал (+/-) p=0,05"

# --- Row 2 (existing SVM / 400 samples result) ---
# "CV" method label becomes "CV, 5" (5-fold cross validation)
$ws.Range("D2").Value = "CV, 5"

# --- Row 3: new result row for SVM run on the last 5000 samples ---
$ws.Range("A3").Value = "SVM (SVC)"
$ws.Range("B3").Value = "Linear, C=1"
$ws.Range("C3").Value = 5000
$ws.Range("D3").Value = "CV, 5"
$ws.Range("E3").Value = 0.73
$ws.Range("F3").Value = 0.02
$ws.Range("G3").Value = 2622.81

# --- Column F width widened to fit the new, longer header text ---
$ws.Columns.Item(6).ColumnWidth = 35.14

# --- Selection moves to G4 ---
$ws.Range("G4").Select()
